$wb = $excel.ActiveWorkbook

$sheetUpdates = @{
    "Citywide Totals" = @(@{ Cell = "I2"; Value = 7276 }, @{ Cell = "J2"; Value = 2379 }, @{ Cell = "I3"; Value = 7487 }, @{ Cell = "J3"; Value = 2467 }, @{ Cell = "J4"; Value = 559 }, @{ Cell = "J5"; Value = 173 }, @{ Cell = "J6"; Value = 3079 }, @{ Cell = "I7"; Value = 26204 }, @{ Cell = "J7"; Value = 8657 });
    "By Neighborhood" = @(@{ Cell = "J2"; Value = 70 }, @{ Cell = "J6"; Value = 85 }, @{ Cell = "J7"; Value = 266 }, @{ Cell = "J8"; Value = 547 }, @{ Cell = "J10"; Value = 53 }, @{ Cell = "J11"; Value = 120 }, @{ Cell = "J13"; Value = 11 }, @{ Cell = "J14"; Value = 32 }, @{ Cell = "J15"; Value = 106 }, @{ Cell = "J18"; Value = 100 }, @{ Cell = "J19"; Value = 276 }, @{ Cell = "J20"; Value = 182 }, @{ Cell = "J23"; Value = 85 }, @{ Cell = "J29"; Value = 485 }, @{ Cell = "J31"; Value = 66 }, @{ Cell = "J33"; Value = 355 }, @{ Cell = "J37"; Value = 293 }, @{ Cell = "J41"; Value = 52 }, @{ Cell = "J42"; Value = 326 }, @{ Cell = "J43"; Value = 82 }, @{ Cell = "J46"; Value = 28 }, @{ Cell = "J47"; Value = 77 }, @{ Cell = "J48"; Value = 85 }, @{ Cell = "J49"; Value = 54 }, @{ Cell = "J50"; Value = 50 }, @{ Cell = "J52"; Value = 215 }, @{ Cell = "J53"; Value = 84 }, @{ Cell = "J54"; Value = 175 }, @{ Cell = "J55"; Value = 103 }, @{ Cell = "J57"; Value = 42 }, @{ Cell = "J60"; Value = 57 }, @{ Cell = "I63"; Value = 208 }, @{ Cell = "J64"; Value = 57 }, @{ Cell = "J65"; Value = 228 }, @{ Cell = "J67"; Value = 314 }, @{ Cell = "J74"; Value = 13 }, @{ Cell = "J76"; Value = 124 }, @{ Cell = "J77"; Value = 66 }, @{ Cell = "J79"; Value = 264 }, @{ Cell = "J83"; Value = 209 }, @{ Cell = "J85"; Value = 403 }, @{ Cell = "J88"; Value = 86 }, @{ Cell = "J91"; Value = 95 }, @{ Cell = "J92"; Value = 28 }, @{ Cell = "J93"; Value = 44 }, @{ Cell = "J95"; Value = 128 }, @{ Cell = "J99"; Value = 118 }, @{ Cell = "J100"; Value = 18 }, @{ Cell = "I101"; Value = 26204 }, @{ Cell = "J101"; Value = 8657 });
    "South Shore" = @(@{ Cell = "J2"; Value = 98 }, @{ Cell = "J3"; Value = 152 }, @{ Cell = "J6"; Value = 118 }, @{ Cell = "J7"; Value = 403 });
    "Little Village" = @(@{ Cell = "J6"; Value = 91 }, @{ Cell = "J7"; Value = 215 });
    "Belmont Cragin" = @(@{ Cell = "J3"; Value = 24 }, @{ Cell = "J7"; Value = 120 });
    "Austin" = @(@{ Cell = "J3"; Value = 178 }, @{ Cell = "J4"; Value = 25 }, @{ Cell = "J6"; Value = 161 }, @{ Cell = "J7"; Value = 547 });
    "Logan Square" = @(@{ Cell = "J6"; Value = 49 }, @{ Cell = "J7"; Value = 84 });
    "Auburn Gresham" = @(@{ Cell = "J2"; Value = 88 }, @{ Cell = "J7"; Value = 266 });
    "Bridgeport" = @(@{ Cell = "J3"; Value = 8 }, @{ Cell = "J7"; Value = 32 });
    "Grand Crossing" = @(@{ Cell = "J2"; Value = 89 }, @{ Cell = "J7"; Value = 293 });
    "Woodlawn" = @(@{ Cell = "J2"; Value = 37 }, @{ Cell = "J7"; Value = 118 });
    "North Lawndale" = @(@{ Cell = "J2"; Value = 66 }, @{ Cell = "J7"; Value = 314 });
    "Gage Park" = @(@{ Cell = "J2"; Value = 26 }, @{ Cell = "J7"; Value = 66 });
    "New City" = @(@{ Cell = "J3"; Value = 63 }, @{ Cell = "J7"; Value = 228 });
    "South Chicago" = @(@{ Cell = "J3"; Value = 74 }, @{ Cell = "J7"; Value = 209 });
    "West Pullman" = @(@{ Cell = "J2"; Value = 47 }, @{ Cell = "J6"; Value = 38 }, @{ Cell = "J7"; Value = 128 });
    "Garfield Park" = @(@{ Cell = "J3"; Value = 103 }, @{ Cell = "J5"; Value = 15 }, @{ Cell = "J7"; Value = 355 });
    "Lincoln Park" = @(@{ Cell = "J3"; Value = 14 }, @{ Cell = "J4"; Value = 4 }, @{ Cell = "J7"; Value = 54 });
    "Loop" = @(@{ Cell = "J4"; Value = 13 }, @{ Cell = "J6"; Value = 84 }, @{ Cell = "J7"; Value = 175 });
    "Englewood" = @(@{ Cell = "J2"; Value = 143 }, @{ Cell = "J6"; Value = 133 }, @{ Cell = "J7"; Value = 485 });
    "Chatham" = @(@{ Cell = "J3"; Value = 76 }, @{ Cell = "J6"; Value = 105 }, @{ Cell = "J7"; Value = 276 });
    "Lake View" = @(@{ Cell = "J3"; Value = 12 }, @{ Cell = "J6"; Value = 41 }, @{ Cell = "J7"; Value = 85 });
    "River North" = @(@{ Cell = "J2"; Value = 17 }, @{ Cell = "J6"; Value = 69 }, @{ Cell = "J7"; Value = 124 });
    "Ashburn" = @(@{ Cell = "J4"; Value = 2 }, @{ Cell = "J7"; Value = 85 });
    "Hermosa" = @(@{ Cell = "J6"; Value = 26 }, @{ Cell = "J7"; Value = 52 });
    "Humboldt Park" = @(@{ Cell = "J3"; Value = 70 }, @{ Cell = "J6"; Value = 167 }, @{ Cell = "J7"; Value = 326 });
    "Boystown" = @(@{ Cell = "J5"; Value = 6 }, @{ Cell = "J6"; Value = 11 });
    "Avondale" = @(@{ Cell = "J4"; Value = 3 }, @{ Cell = "J7"; Value = 53 });
    "Lower West Side" = @(@{ Cell = "J2"; Value = 28 }, @{ Cell = "J6"; Value = 54 }, @{ Cell = "J7"; Value = 103 });
    "Jefferson Park" = @(@{ Cell = "J2"; Value = 9 }, @{ Cell = "J7"; Value = 28 });
    "Douglas" = @(@{ Cell = "J3"; Value = 30 }, @{ Cell = "J6"; Value = 23 }, @{ Cell = "J7"; Value = 85 });
    "Washington Park" = @(@{ Cell = "J3"; Value = 41 }, @{ Cell = "J7"; Value = 95 });
    "Roseland" = @(@{ Cell = "J2"; Value = 75 }, @{ Cell = "J7"; Value = 264 });
    "Near South Side" = @(@{ Cell = "J2"; Value = 17 }, @{ Cell = "J7"; Value = 57 });
    "Chicago Lawn" = @(@{ Cell = "J2"; Value = 57 }, @{ Cell = "J3"; Value = 54 }, @{ Cell = "J7"; Value = 182 });
    "Calumet Heights" = @(@{ Cell = "J6"; Value = 55 }, @{ Cell = "J7"; Value = 100 });
    "West Lawn" = @(@{ Cell = "J3"; Value = 17 }, @{ Cell = "J7"; Value = 44 });
    "Wrigleyville" = @(@{ Cell = "J6"; Value = 8 }, @{ Cell = "J7"; Value = 18 });
    "Kenwood" = @(@{ Cell = "J6"; Value = 33 }, @{ Cell = "J7"; Value = 77 });
    "Brighton Park" = @(@{ Cell = "J3"; Value = 28 }, @{ Cell = "J7"; Value = 106 });
    "Lincoln Square" = @(@{ Cell = "J3"; Value = 17 }, @{ Cell = "J4"; Value = 10 }, @{ Cell = "J7"; Value = 50 });
    "Albany Park" = @(@{ Cell = "J3"; Value = 22 }, @{ Cell = "J6"; Value = 19 }, @{ Cell = "J7"; Value = 70 });
    "West Elsdon" = @(@{ Cell = "J2"; Value = 6 }, @{ Cell = "J7"; Value = 28 });
    "United Center" = @(@{ Cell = "J3"; Value = 26 }, @{ Cell = "J7"; Value = 86 });
    "Mckinley Park" = @(@{ Cell = "J6"; Value = 17 }, @{ Cell = "J7"; Value = 42 });
    "Morgan Park" = @(@{ Cell = "J3"; Value = 15 }, @{ Cell = "J7"; Value = 57 });
    "Hyde Park" = @(@{ Cell = "J3"; Value = 16 }, @{ Cell = "J7"; Value = 82 });
    "Riverdale" = @(@{ Cell = "J2"; Value = 20 }, @{ Cell = "J3"; Value = 21 }, @{ Cell = "J7"; Value = 66 });
    "Printers Row" = @(@{ Cell = "J6"; Value = 7 }, @{ Cell = "J7"; Value = 13 });
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $sheetUpdates[$sheetName]) {
        $ws.Range($u.Cell).Value = $u.Value
    }
}